$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2:A16").Value = "2025-11-25 18:35:56"
